$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(15, "ATL_TB", 47.5, 5.5),
    @(15, "LAC_KC", 45.5, 4.5),
    @(15, "NYJ_JAX", 43.5, 4.5),
    @(15, "BUF_NE", 46.5, -4.5),
    @(15, "ARI_HOU", 45.5, 2.5),
    @(15, "BAL_CIN", 49.5, -2.5),
    @(15, "CLE_CHI", 42.5, 7),
    @(15, "LV_PHI", 44.5, 10.5),
    @(15, "WAS_NYG", 45.5, -4.5),
    @(15, "DET_LA", 48.5, -1.5),
    @(15, "CAR_NO", 43.5, -1.5),
    @(15, "GB_DEN", 44.5, 1.5),
    @(15, "TEN_SF", 45.5, 7.5),
    @(15, "IND_SEA", 44.5, 3.5),
    @(15, "MIN_DAL", 45.5, -2.5),
    @(15, "MIA_PIT", 44.5, 3)
)

$startRow = 194
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Activate()
$ws.Application.Goto($ws.Range("A188"), $true)
$ws.Range("D194").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 188
$excel.ActiveWindow.ScrollColumn = 1
